$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.300082802772522
$ws.Range("B1").Value = 2.249264478683472
$ws.Range("C1").Value = 2.812947511672974
$ws.Range("D1").Value = 3.239881992340088
$ws.Range("E1").Value = 2.07282567024231
